$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44910
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 6500
$ws.Range("P2").Value = 406
$ws.Range("D3").Value = 44698
$ws.Range("J3").Value = 34
$ws.Range("K3").Value = 6000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 6500
$ws.Range("P3").Value = 406
$ws.Range("D4").Value = 44306
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = 6000
$ws.Range("P4").Value = 375
$ws.Range("D5").Value = 44782
$ws.Range("J5").Value = 70
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("D6").Value = 44328
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 6000
$ws.Range("M6").Value = 6000
$ws.Range("P6").Value = 375
$ws.Range("D7").Value = 44943
$ws.Range("J7").Value = 70
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("P7").Value = 438
$ws.Range("D8").Value = 44589
$ws.Range("J8").Value = 52
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 8000
$ws.Range("P8").Value = 500
$ws.Range("D9").Value = 44363
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 5500
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 5750
$ws.Range("P9").Value = 359
$ws.Range("D10").Value = 44355
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 6000
$ws.Range("M10").Value = 6000
$ws.Range("P10").Value = 375
$ws.Range("D11").Value = 44939
$ws.Range("J11").Value = 70
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 438
$ws.Range("D12").Value = 44973
$ws.Range("J12").Value = 90
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 7500
$ws.Range("P12").Value = 469
$ws.Range("D13").Value = 44474
$ws.Range("J13").Value = 52
$ws.Range("K13").Value = 5000
$ws.Range("M13").Value = 5500
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 344
$ws.Range("D14").Value = 44358
$ws.Range("J14").Value = 52
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 6000
$ws.Range("P14").Value = 375
$ws.Range("D15").Value = 44573
$ws.Range("J15").Value = 34
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("P15").Value = 500
$ws.Range("D16").Value = 44582
$ws.Range("J16").Value = 52
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("P16").Value = 438
$ws.Range("D17").Value = 44967
$ws.Range("K17").Value = 7000
$ws.Range("M17").Value = 7000
$ws.Range("P17").Value = 438
$ws.Range("D18").Value = 44938
$ws.Range("D19").Value = 44341
$ws.Range("J19").Value = 51
$ws.Range("K19").Value = 5500
$ws.Range("M19").Value = 5755
$ws.Range("P19").Value = 360
$ws.Range("D20").Value = 44308
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 6000
$ws.Range("P20").Value = 375
$ws.Range("D21").Value = 44715
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 6000
$ws.Range("M21").Value = 5500
$ws.Range("P21").Value = 344
$ws.Range("D22").Value = 44313
$ws.Range("J22").Value = 34
$ws.Range("D23").Value = 44971
$ws.Range("J23").Value = 160
$ws.Range("D24").Value = 44575
$ws.Range("J24").Value = 61
$ws.Range("K24").Value = 8000
$ws.Range("L24").Value = 8000
$ws.Range("M24").Value = 8000
$ws.Range("P24").Value = 500
$ws.Range("D26").Value = 44403
$ws.Range("J26").Value = 43
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 6000
$ws.Range("M26").Value = 6000
$ws.Range("P26").Value = 375
$ws.Range("D27").Value = 44960
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 7000
$ws.Range("M27").Value = 7000
$ws.Range("P27").Value = 438
$ws.Range("D28").Value = 44407
$ws.Range("J28").Value = 45
$ws.Range("K28").Value = 5500
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = 5744
$ws.Range("P28").Value = 359
$ws.Range("D29").Value = 44691
$ws.Range("J29").Value = 61
$ws.Range("L29").Value = 7000
$ws.Range("M29").Value = 6508
$ws.Range("P29").Value = 407
$ws.Range("D30").Value = 44376
$ws.Range("J30").Value = 43
$ws.Range("K30").Value = 4500
$ws.Range("L30").Value = 5000
$ws.Range("M30").Value = 4756
$ws.Range("P30").Value = 297
$ws.Range("D31").Value = 44957
$ws.Range("J31").Value = 70
$ws.Range("K31").Value = 7000
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = 7000
$ws.Range("P31").Value = 438
$ws.Range("D32").Value = 44932
$ws.Range("K32").Value = 6000
$ws.Range("M32").Value = 6500
$ws.Range("P32").Value = 406
$ws.Range("D33").Value = 44350
$ws.Range("J33").Value = 25
$ws.Range("D34").Value = 44455
$ws.Range("J34").Value = 52
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = 5500
$ws.Range("P34").Value = 344
$ws.Range("D35").Value = 44330
$ws.Range("J35").Value = 120
$ws.Range("K35").Value = 6000
$ws.Range("M35").Value = 6000
$ws.Range("P35").Value = 375
$ws.Range("D36").Value = 44946
$ws.Range("J36").Value = 70
$ws.Range("K36").Value = 7000
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = 7000
$ws.Range("P36").Value = 438
$ws.Range("D37").Value = 44438
$ws.Range("J37").Value = 34
$ws.Range("K37").Value = 5000
$ws.Range("M37").Value = 5500
$ws.Range("P37").Value = 344
$ws.Range("D38").Value = 44442
$ws.Range("J38").Value = 25
$ws.Range("K38").Value = 6000
$ws.Range("M38").Value = 6480
$ws.Range("P38").Value = 405
$ws.Range("D39").Value = 44952
$ws.Range("J39").Value = 70
$ws.Range("K39").Value = 7000
$ws.Range("L39").Value = 7000
$ws.Range("M39").Value = 7000
$ws.Range("P39").Value = 438
$ws.Range("D40").Value = 44477
$ws.Range("J40").Value = 25
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = 6000
$ws.Range("P40").Value = 375
$ws.Range("D41").Value = 44966
$ws.Range("K41").Value = 7000
$ws.Range("M41").Value = 7000
$ws.Range("P41").Value = 438
$ws.Range("D42").Value = 44371
$ws.Range("K42").Value = 5500
$ws.Range("L42").Value = 6000
$ws.Range("M42").Value = 5750
$ws.Range("P42").Value = 359
$ws.Range("D43").Value = 44467
$ws.Range("J43").Value = 52
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = 5500
$ws.Range("P43").Value = 344
